$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 940.35297
$ws.Range("I19").Value = 466.66666
$ws.Range("J19").Value = 1198.7273
$ws.Range("K19").Value = 466.66666
$ws.Range("L19").Value = 1198.7273
$ws.Range("M19").Value = -291.66666
$ws.Range("N19").Value = -1548.7273
$ws.Range("H39").Value = 238.42857
$ws.Range("I39").Value = 75
$ws.Range("K39").Value = 225
$ws.Range("M39").Value = 71
$ws.Range("H98").Value = 1679
$ws.Range("I98").Value = 1633.3889
$ws.Range("K98").Value = 1633.3889
$ws.Range("M98").Value = -135.3888999999999
$ws.Range("H122").Value = 1679
$ws.Range("I122").Value = 1633.3889
$ws.Range("K122").Value = 4900.1667
$ws.Range("M122").Value = -2450.1667
$ws.Range("H132").Value = 2240.1833
$ws.Range("I132").Value = 2304.698
$ws.Range("J132").Value = 1751.7142
$ws.Range("K132").Value = 6914.093999999999
$ws.Range("L132").Value = 5255.142599999999
$ws.Range("M132").Value = -4384.093999999999
$ws.Range("N132").Value = -10315.1426
$ws.Range("H135").Value = 1489.6111
$ws.Range("J135").Value = 3251.3333
$ws.Range("L135").Value = 29261.9997
$ws.Range("N135").Value = -34331.9997
$ws.Range("H137").Value = 59971.516
$ws.Range("I137").Value = 83070.17999999999
$ws.Range("K137").Value = 249210.54
$ws.Range("M137").Value = -246660.54

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5143.64
$ws.Range("I32").Value = 4316.5312
$ws.Range("J32").Value = 24994.25
$ws.Range("K32").Value = 4316.5312
$ws.Range("L32").Value = 24994.25
$ws.Range("M32").Value = -4029.5312
$ws.Range("N32").Value = -25568.25
$ws.Range("H61").Value = 7130.914
$ws.Range("I61").Value = 7679.3335
$ws.Range("K61").Value = 7679.3335
$ws.Range("M61").Value = -7467.3335
$ws.Range("H74").Value = 48307.176
$ws.Range("I74").Value = 10456.277
$ws.Range("J74").Value = 184570.4
$ws.Range("K74").Value = 10456.277
$ws.Range("L74").Value = 184570.4
$ws.Range("M74").Value = -9582.277
$ws.Range("N74").Value = -186318.4
$ws.Range("H77").Value = 48307.176
$ws.Range("I77").Value = 10456.277
$ws.Range("J77").Value = 184570.4
$ws.Range("K77").Value = 52281.385
$ws.Range("L77").Value = 922852
$ws.Range("M77").Value = -47913.385
$ws.Range("N77").Value = -931588
$ws.Range("H110").Value = 1030202.44
$ws.Range("I110").Value = 1390072.6
$ws.Range("J110").Value = 2001.7142
$ws.Range("K110").Value = 1390072.6
$ws.Range("L110").Value = 2001.7142
$ws.Range("M110").Value = -1388027.6
$ws.Range("N110").Value = -6091.7142
$ws.Range("H122").Value = 1752004.8
$ws.Range("I122").Value = 1756083.4
$ws.Range("K122").Value = 5268250.199999999
$ws.Range("M122").Value = -5265800.199999999
$ws.Range("H132").Value = 33858.125
$ws.Range("I132").Value = 7134.1055
$ws.Range("J132").Value = 72916.30499999999
$ws.Range("K132").Value = 21402.3165
$ws.Range("L132").Value = 218748.915
$ws.Range("M132").Value = -18872.3165
$ws.Range("N132").Value = -223808.915
$ws.Range("H136").Value = 7130.914
$ws.Range("I136").Value = 7679.3335
$ws.Range("K136").Value = 23038.0005
$ws.Range("M136").Value = -20488.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3440.5
$ws.Range("I20").Value = 2689.4
$ws.Range("J20").Value = 4379.375
$ws.Range("K20").Value = 2689.4
$ws.Range("L20").Value = 4379.375
$ws.Range("M20").Value = -2442.4
$ws.Range("N20").Value = -4873.375

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25566.299
$ws.Range("I31").Value = 8015.706
$ws.Range("J31").Value = 35511.633
$ws.Range("K31").Value = 8015.706
$ws.Range("L31").Value = 35511.633
$ws.Range("M31").Value = -7720.706
$ws.Range("N31").Value = -36101.633
$ws.Range("H34").Value = 25566.299
$ws.Range("I34").Value = 8015.706
$ws.Range("J34").Value = 35511.633
$ws.Range("K34").Value = 8015.706
$ws.Range("L34").Value = 35511.633
$ws.Range("M34").Value = -7813.706
$ws.Range("N34").Value = -35915.633
$ws.Range("H35").Value = 4702.4
$ws.Range("I35").Value = 1946.4286
$ws.Range("J35").Value = 11133
$ws.Range("K35").Value = 1946.4286
$ws.Range("L35").Value = 11133
$ws.Range("M35").Value = -1652.4286
$ws.Range("N35").Value = -11721
$ws.Range("H105").Value = 427.6316
$ws.Range("I105").Value = 368.05554
$ws.Range("K105").Value = 368.05554
$ws.Range("M105").Value = 1378.94446
$ws.Range("H141").Value = 123866.87
$ws.Range("J141").Value = 130857.68
$ws.Range("L141").Value = 130857.68
$ws.Range("N141").Value = -141217.68

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 52750
$ws.Range("J37").Value = 52750
$ws.Range("L37").Value = 158250
$ws.Range("N37").Value = -158474
$ws.Range("H107").Value = 326.66666
$ws.Range("J107").Value = 323.36365
$ws.Range("L107").Value = 970.09095
$ws.Range("N107").Value = -4810.09095
$ws.Range("H133").Value = 3889
$ws.Range("I133").Value = 3889
$ws.Range("K133").Value = 11667
$ws.Range("M133").Value = -6607
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = 0

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10536226
$ws.Range("I70").Value = 15388676
$ws.Range("J70").Value = 22583
$ws.Range("K70").Value = 15388676
$ws.Range("L70").Value = 22583
$ws.Range("M70").Value = -15388406
$ws.Range("N70").Value = -23123
$ws.Range("H73").Value = 10536226
$ws.Range("I73").Value = 15388676
$ws.Range("J73").Value = 22583
$ws.Range("K73").Value = 15388676
$ws.Range("L73").Value = 22583
$ws.Range("M73").Value = -15387740
$ws.Range("N73").Value = -24455

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5696
$ws.Range("I7").Value = 4710.5415
$ws.Range("J7").Value = 6771.0454
$ws.Range("K7").Value = 4710.5415
$ws.Range("L7").Value = 6771.0454
$ws.Range("M7").Value = -4598.5415
$ws.Range("N7").Value = -6995.0454
$ws.Range("H51").Value = 30750
$ws.Range("J51").Value = 30750
$ws.Range("L51").Value = 30750
$ws.Range("N51").Value = -31706
$ws.Range("H53").Value = 24420.6
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 24420.6
$ws.Range("K53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("M53").Value = 24420.6
$ws.Range("N53").Value = -25456.6
$ws.Range("H122").Value = 8562.786
$ws.Range("I122").Value = 7726.7144
$ws.Range("J122").Value = 9398.857
$ws.Range("K122").Value = 23180.1432
$ws.Range("L122").Value = 28196.571
$ws.Range("M122").Value = -20730.1432
$ws.Range("N122").Value = -33096.571
$ws.Range("H126").Value = 5696
$ws.Range("I126").Value = 4710.5415
$ws.Range("J126").Value = 6771.0454
$ws.Range("K126").Value = 14131.6245
$ws.Range("L126").Value = 20313.1362
$ws.Range("M126").Value = -11661.6245
$ws.Range("N126").Value = -25253.1362

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4932.684
$ws.Range("I122").Value = 2859
$ws.Range("J122").Value = 10739
$ws.Range("K122").Value = 8577
$ws.Range("L122").Value = 32217
$ws.Range("M122").Value = -6127
$ws.Range("N122").Value = -37117
$ws.Range("H132").Value = 11493092
$ws.Range("I132").Value = 12992096
$ws.Range("J132").Value = 1000061.44
$ws.Range("K132").Value = 38976288
$ws.Range("L132").Value = 3000184.32
$ws.Range("M132").Value = -38973758
$ws.Range("N132").Value = -3005244.32
